$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 13 (Item #7, RP-E1/RP-E2 connector): add manufacturer/part info and
# refine description to mention it's a stacking header.
$ws.Range("D13").Value = "Adafruit Industries LLC"
$ws.Range("E13").Value = "1979  [1528-1783-ND (DigiKey)]"
$ws.Range("F13").Value = "2x13 (26pin) RP connector (raspberry pi style stacking header, long) BOTTOM MOUNT!"

# Move the active selection to E13, matching the new cursor position left
# behind by this edit.
$ws.Activate()
$ws.Range("E13").Select()
